# Apply "hybrid bold + color" highlighting to quantitative impact metrics
# (percentages, dollar amounts, large numbers) inside specific bullet
# paragraphs of the resume, matching the target commit.
#
# Strategy: for each target paragraph (identified by a unique, stable
# substring), run Find.Execute scoped to that paragraph's Range to locate
# each metric substring in left-to-right order, then apply Bold + the
# highlight color (hex 2C3E50, stored as BGR decimal for Font.Color) to
# just that located sub-range. Word automatically splits the paragraph's
# runs around the (re)formatted sub-range, which reproduces the run-split
# structure seen in the diff.

$d = $word.ActiveDocument

# Highlight color 2C3E50 expressed as the BGR-packed decimal Word's
# Font.Color property expects (0x50 3E 2C -> blue,green,red bytes).
$highlightColor = 5258796

function HighlightMetric($paragraph, $metricText) {
    $r = $paragraph.Range
    $found = $r.Find.Execute($metricText)
    $r.Font.Bold = 1
    $r.Font.Color = $highlightColor
}

function FindParagraphContaining($doc, $uniqueSubstring) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -like "*$uniqueSubstring*") {
            return $p
        }
    }
    return $null
}

# 1. "Discovered systematic race coding errors ... from 23% to 64%"
$p1 = FindParagraphContaining $d "Discovered systematic race coding errors"
HighlightMetric $p1 "23%"
HighlightMetric $p1 "64%"

# 2. "Utilized advanced sampling methods ... ±4.2% to ±2.1% ... 71% to 87% ..."
$p2 = FindParagraphContaining $d "Utilized advanced sampling methods"
HighlightMetric $p2 "±4.2%"
HighlightMetric $p2 "±2.1%"
HighlightMetric $p2 "71%"
HighlightMetric $p2 "87%"

# 3. "Trigonometric algorithm for boundary estimation ... 73.5% ... $4.7M ..."
$p3 = FindParagraphContaining $d "Trigonometric algorithm for boundary estimation"
HighlightMetric $p3 "73.5%"
HighlightMetric $p3 "$4.7M"

# 4. "Built real-time FEC analysis systems ... valued over $2 trillion"
$p4 = FindParagraphContaining $d "Built real-time FEC analysis systems"
HighlightMetric $p4 "$2"

# 5. "Modernized legacy ETL processes ... reducing processing time by 57%"
$p5 = FindParagraphContaining $d "Modernized legacy ETL processes"
HighlightMetric $p5 "57%"

# 6. "Platform impact: Built redistricting system serving 12,847 analysts ..."
$p6 = FindParagraphContaining $d "Platform impact: Built redistricting system serving"
HighlightMetric $p6 "12,847"

# 7. "Revenue generation: Delivered $4.9M additional revenue ..."
$p7 = FindParagraphContaining $d "Revenue generation: Delivered"
HighlightMetric $p7 "$4.9M"

# 8. "23% conversion rate improvement"
$p8 = FindParagraphContaining $d "conversion rate improvement"
HighlightMetric $p8 "23%"

Write-Host "Highlighting applied to 8 paragraphs"
